$d = $word.ActiveDocument

$replacements = @(
    @("20×65=", "21×76="),
    @("99×85=", "54×91="),
    @("59×90=", "38×87="),
    @("15×51=", "61×58="),
    @("89×87=", "43×14="),
    @("96×58=", "36×48="),
    @("22×91=", "25×98="),
    @("42×97=", "43×85="),
    @("25×49=", "24×22="),
    @("27×74=", "25×16="),
    @("29×59=", "74×90="),
    @("94×19=", "27×31="),
    @("21×97=", "13×53="),
    @("65×35=", "91×31="),
    @("63×95=", "65×75="),
    @("73×93=", "84×93="),
    @("30×56=", "68×53="),
    @("66×12=", "92×89="),
    @("87×91=", "31×28="),
    @("94×89=", "17×98="),
    @("55×33=", "32×74="),
    @("33×27=", "89×91="),
    @("51×87=", "55×53="),
    @("29×40=", "25×92="),
    @("71×85=", "12×22=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
